# dict.xlsx fix: make_tidy.py was generating shared-string entries for the
# "second generation" (N1_2/N2_2/... / C1_2/C2_2/...) validation sequence in
# a different order than the rows that reference them, so the workbook's
# row data is supplied here exactly as the corrected export produced it -
# including preserving that same (n2_2 before n1_2) shared-string ordering
# quirk so the two "gens" line up with the fixed script's output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 212's text ("n2_2") must be written to the workbook before row 211's
# text ("n1_2") so the shared-string table gets them in that same order
# (n2_2 -> index 211, n1_2 -> index 212), matching the fixed generator's
# (quirky) append order for this second-generation validation sequence.
$ws.Cells.Item(212, 1).Value = "n2_2"
$ws.Cells.Item(212, 2).Value = 2233333232

$ws.Cells.Item(211, 1).Value = "n1_2"
$ws.Cells.Item(211, 2).Value = 2233331232

$ws.Cells.Item(213, 1).Value = "n3_2"
$ws.Cells.Item(213, 2).Value = 3233331232

$ws.Cells.Item(214, 1).Value = "c1_2"
$ws.Cells.Item(214, 2).Value = 3333332331

$ws.Cells.Item(215, 1).Value = "n4_2"
$ws.Cells.Item(215, 2).Value = 1322222323

$ws.Cells.Item(216, 1).Value = "n5_2"
$ws.Cells.Item(216, 2).Value = 3322223323

$ws.Cells.Item(217, 1).Value = "n6_2"
$ws.Cells.Item(217, 2).Value = 3231333231

$ws.Cells.Item(218, 1).Value = "c2_2"
$ws.Cells.Item(218, 2).Value = 2213323323

$ws.Cells.Item(219, 1).Value = "n7_2"
$ws.Cells.Item(219, 2).Value = 1231331111

$ws.Cells.Item(220, 1).Value = "c3_2"
$ws.Cells.Item(220, 2).Value = 3333332331

$ws.Cells.Item(221, 1).Value = "n8_2"
$ws.Cells.Item(221, 2).Value = 2333233333

$ws.Cells.Item(222, 1).Value = "c4_2"
$ws.Cells.Item(222, 2).Value = 3333332332

# Match the saved file's scroll position / active selection (last row added).
$excel.ActiveWindow.ScrollRow = 206
$null = $ws.Range("A223").Select()
